$wb = $excel.ActiveWorkbook

# 1. Update the "Ready for handoff" -> "Handback transform failed" status text.
# This shared string is used by Overview!E3, Overview!F3, zh-cn!C3, de-de!C3
$newStatus = "Handback transform failed"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus

# 2. Add "Error Detail" values for the handback-transform-failed rows.
$wsZhCn.Range("P3").Value = "Handback file name: gtxms3yt.fpo is different with handoff file name: 075fc7c8-f5b0-46d8-9ed3-5bf09d01b31b.0bebc8a2c12a87e188c85973aa628a8b984643e3.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: gtxms3yt.fpo is different with handoff file name: 075fc7c8-f5b0-46d8-9ed3-5bf09d01b31b.0bebc8a2c12a87e188c85973aa628a8b984643e3.de-de."

# 3. Widen column P (Error Detail) on both language sheets to fit the new text.
# (The saved OOXML "width" attribute is ColumnWidth + 5/6, so back-solve for
# a ColumnWidth that round-trips to exactly width="40".)
$wsZhCn.Columns.Item(16).ColumnWidth = 235/6
$wsDeDe.Columns.Item(16).ColumnWidth = 235/6
